# Auto-generated: update market-data snapshot cells per scheduled-runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6066.5
$ws.Range("I74").Value = 5250
$ws.Range("K74").Value = 5250
$ws.Range("M74").Value = -4314
$ws.Range("H77").Value = 6066.5
$ws.Range("I77").Value = 5250
$ws.Range("K77").Value = 26250
$ws.Range("M77").Value = -21570
$ws.Range("H80").Value = 689.61536
$ws.Range("I80").Value = 706.8570999999999
$ws.Range("J80").Value = 669.5
$ws.Range("K80").Value = 2120.5713
$ws.Range("L80").Value = 2008.5
$ws.Range("M80").Value = -1122.5713
$ws.Range("N80").Value = -4004.5
$ws.Range("H83").Value = 689.61536
$ws.Range("I83").Value = 706.8570999999999
$ws.Range("J83").Value = 669.5
$ws.Range("K83").Value = 6361.7139
$ws.Range("L83").Value = 6025.5
$ws.Range("M83").Value = -1369.7139
$ws.Range("N83").Value = -16009.5
$ws.Range("H86").Value = 1426.8
$ws.Range("I86").Value = 1441
$ws.Range("K86").Value = 1441
$ws.Range("M86").Value = -318
$ws.Range("H89").Value = 1426.8
$ws.Range("I89").Value = 1441
$ws.Range("K89").Value = 7205
$ws.Range("M89").Value = -1589
$ws.Range("H98").Value = 1805.08
$ws.Range("I98").Value = 1581.7727
$ws.Range("K98").Value = 1581.7727
$ws.Range("M98").Value = -83.77269999999999
$ws.Range("H100").Value = 3600.8572
$ws.Range("I100").Value = 4050
$ws.Range("J100").Value = 3002
$ws.Range("K100").Value = 4050
$ws.Range("L100").Value = 3002
$ws.Range("M100").Value = -3509
$ws.Range("N100").Value = -4084
$ws.Range("H106").Value = 1339.8
$ws.Range("I106").Value = 2100
$ws.Range("J106").Value = 833
$ws.Range("K106").Value = 2100
$ws.Range("L106").Value = 833
$ws.Range("M106").Value = -1469
$ws.Range("N106").Value = -2095
$ws.Range("H122").Value = 1805.08
$ws.Range("I122").Value = 1581.7727
$ws.Range("K122").Value = 4745.3181
$ws.Range("M122").Value = -2295.3181
$ws.Range("H134").Value = 82000
$ws.Range("J134").Value = 82000
$ws.Range("L134").Value = 82000
$ws.Range("N134").Value = -92140
$ws.Range("H137").Value = 12503.909
$ws.Range("I137").Value = 3563.5715
$ws.Range("K137").Value = 10690.7145
$ws.Range("M137").Value = -8140.7145
$ws.Range("H141").Value = 5580.381
$ws.Range("I141").Value = 5359.4
$ws.Range("K141").Value = 16078.2
$ws.Range("M141").Value = -10898.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 8333
$ws.Range("J11").Value = 11499.5
$ws.Range("L11").Value = 11499.5
$ws.Range("N11").Value = -11787.5
$ws.Range("H64").Value = 23433.438
$ws.Range("J64").Value = 20329
$ws.Range("L64").Value = 20329
$ws.Range("N64").Value = -20825
$ws.Range("H67").Value = 23433.438
$ws.Range("J67").Value = 20329
$ws.Range("L67").Value = 20329
$ws.Range("N67").Value = -22045
$ws.Range("H98").Value = 172750
$ws.Range("J98").Value = 172750
$ws.Range("L98").Value = 172750
$ws.Range("N98").Value = -178740
$ws.Range("H132").Value = 2712682.8
$ws.Range("I132").Value = 4401.75
$ws.Range("J132").Value = 5898895.5
$ws.Range("K132").Value = 13205.25
$ws.Range("L132").Value = 17696686.5
$ws.Range("M132").Value = -10675.25
$ws.Range("N132").Value = -17701746.5
$ws.Range("H134").Value = 99123.625
$ws.Range("J134").Value = 99123.625
$ws.Range("L134").Value = 99123.625
$ws.Range("N134").Value = -109263.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 630249400
$ws.Range("J33").Value = 630249400
$ws.Range("L33").Value = 630249400
$ws.Range("N33").Value = -630250072
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
$ws.Range("H134").Value = 15399.667
$ws.Range("I134").Value = 9033.786
$ws.Range("J134").Value = 24311.9
$ws.Range("K134").Value = 27101.358
$ws.Range("L134").Value = 72935.70000000001
$ws.Range("M134").Value = -24566.358
$ws.Range("N134").Value = -78005.70000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 15590.459
$ws.Range("J58").Value = 18209.666
$ws.Range("L58").Value = 18209.666
$ws.Range("N58").Value = -18615.666
$ws.Range("H107").Value = 639.6585
$ws.Range("I107").Value = 466.42307
$ws.Range("K107").Value = 466.42307
$ws.Range("M107").Value = 1453.57693
$ws.Range("H136").Value = 15590.459
$ws.Range("J136").Value = 18209.666
$ws.Range("L136").Value = 54628.99800000001
$ws.Range("N136").Value = -59728.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 494.5
$ws.Range("I40").Value = 530.36365
$ws.Range("K40").Value = 2121.4546
$ws.Range("M40").Value = -2052.4546
$ws.Range("H129").Value = 1590.875
$ws.Range("I129").Value = 1097.8334
$ws.Range("J129").Value = 3070
$ws.Range("K129").Value = 3293.5002
$ws.Range("L129").Value = 9210
$ws.Range("M129").Value = 1706.4998
$ws.Range("N129").Value = -19210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 99998.5
$ws.Range("J10").Value = 99998.5
$ws.Range("L10").Value = 99998.5
$ws.Range("N10").Value = -100336.5
$ws.Range("H97").Value = 1422.55
$ws.Range("I97").Value = 1698.8889
$ws.Range("J97").Value = 1196.4546
$ws.Range("K97").Value = 1698.8889
$ws.Range("L97").Value = 1196.4546
$ws.Range("M97").Value = -1202.8889
$ws.Range("N97").Value = -2188.4546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 977.65625
$ws.Range("I16").Value = 977.65625
$ws.Range("K16").Value = 977.65625
$ws.Range("M16").Value = -807.65625
$ws.Range("H40").Value = 4852.0356
$ws.Range("I40").Value = 4395.1177
$ws.Range("K40").Value = 4395.1177
$ws.Range("M40").Value = -4259.1177
$ws.Range("H46").Value = 402979.6
$ws.Range("J46").Value = 3087.9412
$ws.Range("L46").Value = 3087.9412
$ws.Range("N46").Value = -3463.9412
$ws.Range("H93").Value = 13276
$ws.Range("I93").Value = 14216.667
$ws.Range("K93").Value = 14216.667
$ws.Range("M93").Value = -12968.667
$ws.Range("H100").Value = 3566.6
$ws.Range("I100").Value = 3155.4443
$ws.Range("J100").Value = 4183.3335
$ws.Range("K100").Value = 3155.4443
$ws.Range("M100").Value = -2614.4443
$ws.Range("N100").Value = -5265.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1012.8182
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 1012.8182
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H96").Value = 1998.7778
$ws.Range("I96").Value = 1597
$ws.Range("K96").Value = 1597
$ws.Range("M96").Value = -224
$ws.Range("H100").Value = 849
$ws.Range("I100").Value = 710.4
$ws.Range("J100").Value = 1022.25
$ws.Range("K100").Value = 1420.8
$ws.Range("L100").Value = 2044.5
$ws.Range("M100").Value = -879.8
$ws.Range("N100").Value = -3126.5
